{"js": "// Update the signature location/date, and tidy the \"Kepada Yth IBI Cabang\n// Banyuasin\" line (removing the stray leftover _GoBack bookmark there by\n// merging it back into its text run), per:\n//   update(surat pernyataan, surat permohonan sipb ibi): lokasi tanda tangan\n//   dan tanggal\n\nconst body = context.document.body;\n\n// --- 1) \"IBI Cabang Banyuasin\" line -----------------------------------\n// Was two runs (\"IBI Cabang Banyuasin\" + bookmark _GoBack + \" \"); collapse\n// back into a single run with the same visible text \"IBI Cabang Banyuasin \".\nconst greeting = body.search(\"IBI Cabang Banyuasin \", { matchCase: true });\ngreeting.load(\"items\");\nawait context.sync();\n\nif (greeting.items.length > 0) {\n  greeting.items[0].insertText(\"IBI Cabang Banyuasin \", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// The old _GoBack bookmark (if present) lived on that line; drop it so it\n// doesn't linger on the merged run.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- 2) Signature location/date line ----------------------------------\n// \"Pangkalan Balai, 10 Mei 2024\" -> \"Tanah Mas, 10 Juli 2024\"\nconst dateLine = body.search(\"Pangkalan Balai, 10 Mei 2024\", { matchCase: true });\ndateLine.load(\"items\");\nawait context.sync();\n\nif (dateLine.items.length > 0) {\n  dateLine.items[0].insertText(\"Tanah Mas, 10 Juli 2024\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Word leaves a _GoBack bookmark marking the last edit location; here that\n// lands right after \"...10 Ju\", before \"li 2024\".\nconst editSpot = body.search(\"Ju\", { matchCase: true });\neditSpot.load(\"items\");\nawait context.sync();\n\nif (editSpot.items.length > 0) {\n  const afterJu = editSpot.items[0].getRange(Word.RangeLocation.after);\n  afterJu.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Update the signature location/date, and tidy the \"Kepada Yth IBI Cabang\n# Banyuasin\" line (removing the stray leftover _GoBack bookmark there by\n# merging it back into its text run), per:\n#   update(surat pernyataan, surat permohonan sipb ibi): lokasi tanda tangan\n#   dan tanggal\n\n$d = $word.ActiveDocument\n\n# --- 1) \"IBI Cabang Banyuasin\" line -------------------------------------\n# Was two runs (\"IBI Cabang Banyuasin\" + bookmark _GoBack + \" \"); collapse\n# back into a single run with the same visible text \"IBI Cabang Banyuasin \".\n$greeting = $d.Content\n$null = $greeting.Find.Execute(\"IBI Cabang Banyuasin \", $false, $false, $false, $false, $false, $true, 1, $false, \"IBI Cabang Banyuasin \", 2)\n\n# The old _GoBack bookmark (if present) lived on that line; drop it so it\n# doesn't linger on the merged run.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# --- 2) Signature location/date line ------------------------------------\n# \"Pangkalan Balai, 10 Mei 2024\" -> \"Tanah Mas, 10 Juli 2024\"\n$dateLine = $d.Content\n$null = $dateLine.Find.Execute(\"Pangkalan Balai, 10 Mei 2024\", $false, $false, $false, $false, $false, $true, 1, $false, \"Tanah Mas, 10 Juli 2024\", 2)\n\n# Word leaves a _GoBack bookmark marking the last edit location; here that\n# lands right after \"...10 Ju\", before \"li 2024\".\n$editSpot = $d.Content\n$found = $editSpot.Find.Execute(\"Ju\", $true)\nif ($found) {\n    $editSpot.Collapse(0)\n    $d.Bookmarks.Add(\"_GoBack\", $editSpot)\n}\n"}
